# Refresh cryptocurrency Price and Volume(1h) columns with the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.241.27"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "'2.477.04"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'577.49"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "'146.80"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'2.474.89"
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("D10").Value = "'0.112"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").Value = "'5.29"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").Value = "'29.09"
$ws.Range("E14").Value = "  +8.88%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "'2.924.14"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "'63.162.85"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("D18").Value = "'2.469.92"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").Value = "'8.12"
$ws.Range("D20").Value = "'11.07"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").Value = "'330.76"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").Value = "'2.23"
$ws.Range("E22").Value = "  +10.35%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'66.45"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").Value = "'671.24"
$ws.Range("E26").Value = "  +9.26%  "
$ws.Range("D27").Value = "'9.43"
$ws.Range("E27").Value = "  +12.98%  "
$ws.Range("D28").Value = "'0.0₃0998"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("D29").Value = "'2.600.57"
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("D30").Value = "'0.995"
$ws.Range("E30").Value = "  +874.99%  "
$ws.Range("E31").Value = "  +4.34%  "
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").Value = "'1.87"
$ws.Range("E33").Value = "  +2.97%  "
$ws.Range("E34").Value = "  -2.14%  "
$ws.Range("E35").Value = "  +5.24%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Value = "'4.80"
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").Value = "'153.56"
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("D41").Value = "'18.77"
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D42").Value = "'2.73"
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'0.0₆0310"
$ws.Range("E45").Value = "  +12.25%  "
$ws.Range("D46").Value = "'15.17"
$ws.Range("E46").Value = "  +27.63%  "
$ws.Range("D47").Value = "'149.15"
$ws.Range("E47").Value = "  +4.07%  "
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").Value = "'20.90"
$ws.Range("E49").Value = "  +3.71%  "
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").Value = "'0.0516"
$ws.Range("E51").Value = "  +0.60%  "
